# Weekly update: a new price observation was recorded for
# "Terminal Hortofrutícola Agro Chillán" - Mango, inserted as the new
# row 168. All the historical rows that used to occupy rows 168-196
# shift down by one (to 169-197); Excel's native row-insert does this
# for us automatically, including updating the used-range dimension.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 168 - pushes existing rows 168:196 down to 169:197.
$ws.Rows("168:168").Insert()

# Populate the newly inserted row with this week's data point.
$ws.Range("A168").Value = 7
$ws.Range("B168").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C168").Value = "Ñuble"
$ws.Range("D168").Value = 45209
$ws.Range("E168").Value = 16
$ws.Range("F168").Value = "Fruta"
$ws.Range("G168").Value = 100108
$ws.Range("H168").Value = "Tropicales y subtropicales"
$ws.Range("I168").Value = 100108002
$ws.Range("J168").Value = "Mango"
$ws.Range("K168").Value = "Sin especificar"
$ws.Range("L168").Value = "Primera"
$ws.Range("M168").Value = 80
$ws.Range("N168").Value = 10000
$ws.Range("O168").Value = 10000
$ws.Range("P168").Value = 10000
$ws.Range("Q168").Value = "$/bandeja 4 kilos"
$ws.Range("R168").Value = "Brasil"
$ws.Range("S168").Value = 2500
$ws.Range("T168").Value = 4
